# Update the "cryptos" price list with the latest scraped values.
# Note: several Price values look numeric (e.g. "1.001", "307.81") but must
# stay as literal text (leading zeros / trailing zeros / thousand-dot style
# must be preserved exactly as scraped), so those are written with a
# leading apostrophe to force Excel to store them as text instead of
# silently parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.285.00"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.908.48"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'307.81"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").Value = "'0.3787"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("D9").Value = "'0.07267"
$ws.Range("D10").Value = "'21.33"
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("D11").Value = "'0.9017"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'0.08224"
$ws.Range("E12").Value = "  +9.90%  "
$ws.Range("D13").Value = "1.908.52"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "'95.30"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'5.295"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'0.000008620"
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("D18").Value = "'14.50"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.336.69"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'5.076"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "2.153.06"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("E23").Value = "  +3.32%  "
$ws.Range("D24").Value = "'6.469"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'2.309"
$ws.Range("E25").Value = "  +10.61%  "
$ws.Range("D26").Value = "'146.27"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'1.748"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "'115.05"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  +6.48%  "
$ws.Range("D31").Value = "'4.814"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("D32").Value = "'0.09233"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'0.8079"
$ws.Range("E33").Value = "  +8.00%  "
$ws.Range("D34").Value = "'0.05059"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("E35").Value = "  +8.05%  "
$ws.Range("D36").Value = "'2.972"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'3.370"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "'0.5746"
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").Value = "'0.01984"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'1.079"
$ws.Range("E41").Value = "  +0.61%  "

# Rows 42/43 swapped rank order: Quant now ranks above Aptos.
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'119.68"
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.999"
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").Value = "'6.631"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "'0.1519"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D47").Value = "'10.26"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("D50").Value = "'37.69"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'63.89"
$ws.Range("E51").Value = "  +1.52%  "
